# Weekly data update: insert a new weekly record row for
# "Terminal Hortofrutícola Agro Chillán" / Zanahoria (Ñuble), pushing the
# existing rows 146:280 down to 147:281, and populate the freshly inserted
# row 146 with the new week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 146 (shifts rows 146:280 -> 147:281).
$ws.Rows.Item(146).Insert()

# Populate the new row 146 with this week's observation. The descriptive /
# categorical columns (market, region, product, quality, unit, origin,
# classification) repeat the same values used throughout this sub-sheet.
$ws.Cells.Item(146, 1).Value = 7
$ws.Cells.Item(146, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(146, 3).Value = "Ñuble"
$ws.Cells.Item(146, 4).Value = 44669
$ws.Cells.Item(146, 5).Value = 16
$ws.Cells.Item(146, 6).Value = 100114013
$ws.Cells.Item(146, 7).Value = "Zanahoria"
$ws.Cells.Item(146, 8).Value = "Sin especificar"
$ws.Cells.Item(146, 9).Value = "Primera"
$ws.Cells.Item(146, 10).Value = 100
$ws.Cells.Item(146, 11).Value = 6000
$ws.Cells.Item(146, 12).Value = 6500
$ws.Cells.Item(146, 13).Value = 6250
$ws.Cells.Item(146, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(146, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(146, 16).Value = 312
$ws.Cells.Item(146, 17).Value = 20
$ws.Cells.Item(146, 18).Value = "Hortaliza"
